$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# "Recorded By" column (G) lists the recorder(s) for each attendance session.
# Reorder the combined "System, <email>" entries to "<email>, System" to
# match the updated attendance export ordering.
$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
